$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 28535.2
$ws.Range("I54").Value = 28535.2
$ws.Range("K54").Value = 28535.2
$ws.Range("M54").Value = -28049.2

$ws.Range("H62").Value = 5857.8335
$ws.Range("I62").Value = 2451
$ws.Range("J62").Value = 7561.25
$ws.Range("K62").Value = 2451
$ws.Range("L62").Value = 7561.25
$ws.Range("M62").Value = -1827
$ws.Range("N62").Value = -8809.25

$ws.Range("H64").Value = 7321.2
$ws.Range("J64").Value = 7558.4644
$ws.Range("L64").Value = 7558.4644
$ws.Range("N64").Value = -8054.4644

$ws.Range("H65").Value = 5857.8335
$ws.Range("I65").Value = 2451
$ws.Range("J65").Value = 7561.25
$ws.Range("K65").Value = 12255
$ws.Range("L65").Value = 37806.25
$ws.Range("M65").Value = -9135
$ws.Range("N65").Value = -44046.25

$ws.Range("H67").Value = 7321.2
$ws.Range("J67").Value = 7558.4644
$ws.Range("L67").Value = 7558.4644
$ws.Range("N67").Value = -9274.464400000001

$ws.Range("H76").Value = 6826
$ws.Range("I76").Value = 6624
$ws.Range("J76").Value = 7143.4287
$ws.Range("K76").Value = 6624
$ws.Range("L76").Value = 7143.4287
$ws.Range("M76").Value = -6309
$ws.Range("N76").Value = -7773.4287

$ws.Range("H79").Value = 6826
$ws.Range("I79").Value = 6624
$ws.Range("J79").Value = 7143.4287
$ws.Range("K79").Value = 6624
$ws.Range("L79").Value = 7143.4287
$ws.Range("M79").Value = -5532
$ws.Range("N79").Value = -9327.4287

$ws.Range("H98").Value = 1043.4615
$ws.Range("I98").Value = 985.2
$ws.Range("K98").Value = 985.2
$ws.Range("M98").Value = 512.8

$ws.Range("H111").Value = 8550678
$ws.Range("J111").Value = 3349.625
$ws.Range("L111").Value = 10048.875
$ws.Range("N111").Value = -16182.875

$ws.Range("H113").Value = 4807.5386
$ws.Range("J113").Value = 5367
$ws.Range("L113").Value = 5367
$ws.Range("N113").Value = -11875

$ws.Range("H122").Value = 1043.4615
$ws.Range("I122").Value = 985.2
$ws.Range("K122").Value = 2955.6
$ws.Range("M122").Value = -505.6000000000004

$ws.Range("H131").Value = 4703.8647
$ws.Range("I131").Value = 1215.8823
$ws.Range("K131").Value = 3647.6469
$ws.Range("M131").Value = 1392.3531

$ws.Range("H137").Value = 56410.305
$ws.Range("I137").Value = 82880.87
$ws.Range("J137").Value = 3469.182
$ws.Range("K137").Value = 248642.61
$ws.Range("L137").Value = 10407.546
$ws.Range("M137").Value = -246092.61
$ws.Range("N137").Value = -15507.546

$ws.Range("H138").Value = 3323.9836
$ws.Range("I138").Value = 2470.7334
$ws.Range("J138").Value = 3602.2173
$ws.Range("K138").Value = 7412.2002
$ws.Range("L138").Value = 10806.6519
$ws.Range("M138").Value = -2272.2002
$ws.Range("N138").Value = -21086.6519

$ws.Range("H141").Value = 5284.074
$ws.Range("I141").Value = 5284.074
$ws.Range("K141").Value = 15852.222
$ws.Range("M141").Value = -10672.222

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1274.39
$ws.Range("I32").Value = 1274.39
$ws.Range("K32").Value = 1274.39
$ws.Range("M32").Value = -987.3900000000001

$ws.Range("H55").Value = 50000
$ws.Range("J55").Value = 50000
$ws.Range("L55").Value = 50000
$ws.Range("N55").Value = -50630

$ws.Range("H60").Value = 8999.666999999999
$ws.Range("I60").Value = 8999.666999999999
$ws.Range("K60").Value = 8999.666999999999
$ws.Range("M60").Value = -8266.666999999999

$ws.Range("H74").Value = 14069.827
$ws.Range("I74").Value = 3797.5
$ws.Range("J74").Value = 63377
$ws.Range("K74").Value = 3797.5
$ws.Range("L74").Value = 63377
$ws.Range("M74").Value = -2923.5
$ws.Range("N74").Value = -65125

$ws.Range("H77").Value = 14069.827
$ws.Range("I77").Value = 3797.5
$ws.Range("J77").Value = 63377
$ws.Range("K77").Value = 18987.5
$ws.Range("L77").Value = 316885
$ws.Range("M77").Value = -14619.5
$ws.Range("N77").Value = -325621

$ws.Range("H110").Value = 1463466.1
$ws.Range("I110").Value = 1635285.8
$ws.Range("K110").Value = 1635285.8
$ws.Range("M110").Value = -1633240.8

$ws.Range("H122").Value = 3677593.8
$ws.Range("I122").Value = 5265458
$ws.Range("K122").Value = 15796374
$ws.Range("M122").Value = -15793924

$ws.Range("H132").Value = 28176.871
$ws.Range("I132").Value = 6609.25
$ws.Range("K132").Value = 19827.75
$ws.Range("M132").Value = -17297.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3031.1738
$ws.Range("I20").Value = 2498.125
$ws.Range("J20").Value = 4249.5713
$ws.Range("K20").Value = 2498.125
$ws.Range("L20").Value = 4249.5713
$ws.Range("M20").Value = -2251.125
$ws.Range("N20").Value = -4743.5713

$ws.Range("H86").Value = 12750577
$ws.Range("I86").Value = 27085072
$ws.Range("K86").Value = 27085072
$ws.Range("M86").Value = -27083949

$ws.Range("H89").Value = 12750577
$ws.Range("I89").Value = 27085072
$ws.Range("K89").Value = 135425360
$ws.Range("M89").Value = -135419744

$ws.Range("H94").Value = 2225038.8
$ws.Range("I94").Value = 2565117
$ws.Range("K94").Value = 2565117
$ws.Range("M94").Value = -2564666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1862
$ws.Range("I16").Value = 1547.8
$ws.Range("J16").Value = 2647.5
$ws.Range("K16").Value = 1547.8
$ws.Range("L16").Value = 2647.5
$ws.Range("M16").Value = -1260.8
$ws.Range("N16").Value = -3221.5

$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H62").Value = 6833
$ws.Range("I62").Value = 6833
$ws.Range("K62").Value = 6833
$ws.Range("M62").Value = -6209

$ws.Range("H65").Value = 6833
$ws.Range("I65").Value = 6833
$ws.Range("K65").Value = 34165
$ws.Range("M65").Value = -31045

$ws.Range("H113").Value = 1862
$ws.Range("I113").Value = 1547.8
$ws.Range("J113").Value = 2647.5
$ws.Range("K113").Value = 1547.8
$ws.Range("L113").Value = 2647.5
$ws.Range("M113").Value = 622.2
$ws.Range("N113").Value = -6987.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 60166.668
$ws.Range("J37").Value = 60166.668
$ws.Range("L37").Value = 180500.004
$ws.Range("N37").Value = -180724.004

$ws.Range("H132").Value = 1494.2084
$ws.Range("J132").Value = 1377.5333
$ws.Range("L132").Value = 12397.7997
$ws.Range("N132").Value = -17457.7997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws.Range("H70").Value = 10009673
$ws.Range("I70").Value = 22226084
$ws.Range("J70").Value = 14428.363
$ws.Range("K70").Value = 22226084
$ws.Range("L70").Value = 14428.363
$ws.Range("M70").Value = -22225814
$ws.Range("N70").Value = -14968.363

$ws.Range("H73").Value = 10009673
$ws.Range("I73").Value = 22226084
$ws.Range("J73").Value = 14428.363
$ws.Range("K73").Value = 22226084
$ws.Range("L73").Value = 14428.363
$ws.Range("M73").Value = -22225148
$ws.Range("N73").Value = -16300.363

$ws.Range("H80").Value = 2099976.5
$ws.Range("I80").Value = 3535670.5
$ws.Range("J80").Value = 425000
$ws.Range("K80").Value = 3535670.5
$ws.Range("L80").Value = 425000
$ws.Range("M80").Value = -3534672.5
$ws.Range("N80").Value = -426996

$ws.Range("H83").Value = 2099976.5
$ws.Range("I83").Value = 3535670.5
$ws.Range("J83").Value = 425000
$ws.Range("K83").Value = 17678352.5
$ws.Range("L83").Value = 2125000
$ws.Range("M83").Value = -17673360.5
$ws.Range("N83").Value = -2134984

$ws.Range("H107").Value = 7726.2144
$ws.Range("I107").Value = 11533
$ws.Range("K107").Value = 11533
$ws.Range("M107").Value = -9613

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4504.4893
$ws.Range("I40").Value = 2597.2354
$ws.Range("J40").Value = 9492.691999999999
$ws.Range("K40").Value = 2597.2354
$ws.Range("L40").Value = 9492.691999999999
$ws.Range("M40").Value = -2461.2354
$ws.Range("N40").Value = -9764.691999999999

$ws.Range("H46").Value = 4841.826
$ws.Range("I46").Value = 1128
$ws.Range("K46").Value = 1128
$ws.Range("M46").Value = -940

$ws.Range("H68").Value = 1750.5555
$ws.Range("I68").Value = 888.3333
$ws.Range("K68").Value = 888.3333
$ws.Range("M68").Value = -139.3333

$ws.Range("H71").Value = 1750.5555
$ws.Range("I71").Value = 888.3333
$ws.Range("K71").Value = 4441.6665
$ws.Range("M71").Value = -697.6665000000003

$ws.Range("H82").Value = 2417743.5
$ws.Range("I82").Value = 3706291.2
$ws.Range("K82").Value = 3706291.2
$ws.Range("M82").Value = -3705930.2

$ws.Range("H85").Value = 2417743.5
$ws.Range("I85").Value = 3706291.2
$ws.Range("K85").Value = 3706291.2
$ws.Range("M85").Value = -3705043.2
